# Weekly update: insert two new price observation rows into the daily
# logic subset sheet ("Hortaliza, Femacal de La Calera - Poroto granado").
# One new row is inserted right after the header (most recent record),
# and a second new row is inserted further down the list; all other
# existing rows shift down to make room, preserving their original data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row right after the header (becomes row 45) ---
$ws.Rows.Item(45).Insert()

$ws.Range("A45").Value = 3
$ws.Range("B45").Value = "Femacal de La Calera"
$ws.Range("C45").Value = "Coquimbo"
$ws.Range("D45").Value = 44614
$ws.Range("E45").Value = 5
$ws.Range("F45").Value = 100112030
$ws.Range("G45").Value = "Poroto granado"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 73
$ws.Range("K45").Value = 24000
$ws.Range("L45").Value = 25000
$ws.Range("M45").Value = 24521
$ws.Range("N45").Value = "`$/saco 25 kilos"
$ws.Range("O45").Value = "Provincia de Petorca"
$ws.Range("P45").Value = 981
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"

# --- Insert a second new row further down the list (becomes row 142) ---
$ws.Rows.Item(142).Insert()

$ws.Range("A142").Value = 3
$ws.Range("B142").Value = "Femacal de La Calera"
$ws.Range("C142").Value = "Coquimbo"
$ws.Range("D142").Value = 44615
$ws.Range("E142").Value = 5
$ws.Range("F142").Value = 100112030
$ws.Range("G142").Value = "Poroto granado"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 65
$ws.Range("K142").Value = 24000
$ws.Range("L142").Value = 25000
$ws.Range("M142").Value = 24462
$ws.Range("N142").Value = "`$/saco 25 kilos"
$ws.Range("O142").Value = "Provincia de Limarí"
$ws.Range("P142").Value = 978
$ws.Range("Q142").Value = 25
$ws.Range("R142").Value = "Hortaliza"
